$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 523 (shifts existing rows 523:641 down to 524:642)
$ws.Rows("523:523").Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A523").Value = 3
$ws.Range("B523").Value = "Femacal de La Calera"
$ws.Range("C523").Value = "Coquimbo"
$ws.Range("D523").Value = 45244
$ws.Range("E523").Value = 5
$ws.Range("F523").Value = 100112012
$ws.Range("G523").Value = "Espinaca"
$ws.Range("H523").Value = "Sin especificar"
$ws.Range("I523").Value = "Primera"
$ws.Range("J523").Value = 80
$ws.Range("K523").Value = 4000
$ws.Range("L523").Value = 4000
$ws.Range("M523").Value = 4000
$ws.Range("N523").Value = '$/docena de atados (3 kilos)'
$ws.Range("O523").Value = "Provincia de Quillota"
$ws.Range("P523").Value = 1333
$ws.Range("Q523").Value = 3
$ws.Range("R523").Value = "Hortaliza"
